$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the cell values below
# (which live on the protected sheet) can be updated, then restore protection.
$ws.Unprotect()

# Update the confidentiality footer text: "as of" date moves from 2021-05-20 to 2021-05-21
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Refreshed Weight (D) and Percent Change (E) figures for rows 2-12
$ws.Range("D2").Value = 0.03104740954994176
$ws.Range("E2").Value = -0.0006823609689524668

$ws.Range("D3").Value = 0.0237518473374715
$ws.Range("E3").Value = 0.005887894488930678

$ws.Range("D4").Value = 0.05245706446724258
$ws.Range("E4").Value = 0.0002314279102058059

$ws.Range("D5").Value = 0.1392425597988659
$ws.Range("E5").Value = -0.0003231017770598532

$ws.Range("D6").Value = 0.03142648790037845
$ws.Range("E6").Value = 0.002832861189801639

$ws.Range("D7").Value = 0.1164197817246625
$ws.Range("E7").Value = -0.002986237340950271

$ws.Range("D8").Value = 0.1016041857867161
$ws.Range("E8").Value = 0.004812141402924297

$ws.Range("D9").Value = 0.02934244962215351
$ws.Range("E9").Value = 0.002689284236656819

$ws.Range("D10").Value = 0.1264943393545099
$ws.Range("E10").Value = 0.008280015054572765

$ws.Range("D11").Value = 0.2455178001982157
$ws.Range("E11").Value = -0.005065856129685908

$ws.Range("D12").Value = 0.1026960742598421
$ws.Range("E12").Value = -0.002725856697819329

# Row 13 (Total): only Percent Change (E) changes
$ws.Range("E13").Value = -0.0000812905850386958

# Restore sheet protection to match the original state
$ws.Protect()
